$d = $word.ActiveDocument

# Sequential Find & Replace, in document order, matching the exact text of each
# run. Each "old" value is unique at the time it is searched for (forward search
# from the start of the document with wrap disabled), so this reproduces the diff
# exactly even though some replacement targets coincide with other cells' values.
$d.Content.Find.Execute("2024-05-25 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-26 Sunday", 2) | Out-Null
$d.Content.Find.Execute("98÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷6=", 2) | Out-Null
$d.Content.Find.Execute("69÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=", 2) | Out-Null
$d.Content.Find.Execute("43÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷7=", 2) | Out-Null
$d.Content.Find.Execute("66÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷8=", 2) | Out-Null
$d.Content.Find.Execute("33÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=", 2) | Out-Null
$d.Content.Find.Execute("69÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷4=", 2) | Out-Null
$d.Content.Find.Execute("72÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=", 2) | Out-Null
$d.Content.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷2=", 2) | Out-Null
$d.Content.Find.Execute("98÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷7=", 2) | Out-Null
$d.Content.Find.Execute("67÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷6=", 2) | Out-Null
$d.Content.Find.Execute("48÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=", 2) | Out-Null
$d.Content.Find.Execute("66÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=", 2) | Out-Null
$d.Content.Find.Execute("94÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷3=", 2) | Out-Null
$d.Content.Find.Execute("85÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=", 2) | Out-Null
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=", 2) | Out-Null
$d.Content.Find.Execute("37÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷3=", 2) | Out-Null
$d.Content.Find.Execute("99÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷9=", 2) | Out-Null
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷6=", 2) | Out-Null
$d.Content.Find.Execute("99÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷4=", 2) | Out-Null
$d.Content.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=", 2) | Out-Null
$d.Content.Find.Execute("96÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=", 2) | Out-Null
$d.Content.Find.Execute("50÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷3=", 2) | Out-Null
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷7=", 2) | Out-Null
$d.Content.Find.Execute("52÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷7=", 2) | Out-Null
$d.Content.Find.Execute("72÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷7=", 2) | Out-Null
